$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.472738
$ws.Range("H2").Value = 1.418214
$ws.Range("I2").Value = 0.0327564895931267
$ws.Range("J2").Value = 0.03397138804734427
$ws.Range("M2").Value = 16.14072933333334
$ws.Range("N2").Value = 48.42218800000001
$ws.Range("O2").Value = 0.03423048004954622
$ws.Range("P2").Value = 0.03634868370049611
$ws.Range("Q2").Value = 7.630336103581335
$ws.Range("R2").Value = 68.67302493223201
$ws.Range("S2").Value = 0.001121270363510692
$ws.Range("T2").Value = 0.001234815238999731

$ws.Range("G3").Value = 0.472738
$ws.Range("H3").Value = 1.418214
$ws.Range("I3").Value = 0.0327564895931267
$ws.Range("J3").Value = 0.03397138804734427
$ws.Range("O3").Value = 0.1719151703242873
$ws.Range("P3").Value = 0.1825533892714798
$ws.Range("Q3").Value = 38.32171003678734
$ws.Range("R3").Value = 344.895390331086
$ws.Range("S3").Value = 0.005631337487628122
$ws.Range("T3").Value = 0.006201592026299337

$ws.Range("G4").Value = 0.472738
$ws.Range("H4").Value = 1.418214
$ws.Range("I4").Value = 0.0327564895931267
$ws.Range("J4").Value = 0.03397138804734427
$ws.Range("M4").Value = 168.70371
$ws.Range("N4").Value = 506.11113
$ws.Range("O4").Value = 0.3577786889414888
$ws.Range("P4").Value = 0.3799182594076638
$ws.Range("Q4").Value = 79.75265445798001
$ws.Range("R4").Value = 717.77389012182
$ws.Range("S4").Value = 0.01171957390095439
$ws.Range("T4").Value = 0.01290635061660935

$ws.Range("G5").Value = 0.472738
$ws.Range("H5").Value = 1.418214
$ws.Range("I5").Value = 0.0327564895931267
$ws.Range("J5").Value = 0.03397138804734427
$ws.Range("M5").Value = 82.43477250000001
$ws.Range("N5").Value = 164.869545
$ws.Range("O5").Value = 0.1748236883957081
$ws.Range("P5").Value = 0.1237612588479007
$ws.Range("Q5").Value = 38.970049482105
$ws.Range("R5").Value = 233.82029689263
$ws.Range("S5").Value = 0.005726610329566036
$ws.Range("T5").Value = 0.004204341749549854

$ws.Range("G6").Value = 0.472738
$ws.Range("H6").Value = 1.418214
$ws.Range("I6").Value = 0.0327564895931267
$ws.Range("J6").Value = 0.03397138804734427
$ws.Range("M6").Value = 123.1883796666667
$ws.Range("N6").Value = 369.565139
$ws.Range("O6").Value = 0.2612519722889696
$ws.Range("P6").Value = 0.2774184087724594
$ws.Range("Q6").Value = 58.23582822686067
$ws.Range("R6").Value = 524.122454041746
$ws.Range("S6").Value = 0.008557697511467456
$ws.Range("T6").Value = 0.009424288415885995

$ws.Range("I7").Value = 0.822180234441485
$ws.Range("J7").Value = 0.8526739017519405
$ws.Range("M7").Value = 16.14072933333334
$ws.Range("N7").Value = 48.42218800000001
$ws.Range("O7").Value = 0.03423048004954622
$ws.Range("P7").Value = 0.03634868370049611
$ws.Range("Q7").Value = 191.5196531873245
$ws.Range("R7").Value = 1723.67687868592
$ws.Range("S7").Value = 0.02814362411218049
$ws.Range("T7").Value = 0.03099357395444918

$ws.Range("I8").Value = 0.822180234441485
$ws.Range("J8").Value = 0.8526739017519405
$ws.Range("O8").Value = 0.1719151703242873
$ws.Range("P8").Value = 0.1825533892714798
$ws.Range("S8").Value = 0.1413452550412704
$ws.Range("T8").Value = 0.1556585107081536

$ws.Range("I9").Value = 0.822180234441485
$ws.Range("J9").Value = 0.8526739017519405
$ws.Range("M9").Value = 168.70371
$ws.Range("N9").Value = 506.11113
$ws.Range("O9").Value = 0.3577786889414888
$ws.Range("P9").Value = 0.3799182594076638
$ws.Range("Q9").Value = 2001.7729907588
$ws.Range("R9").Value = 18015.9569168292
$ws.Range("S9").Value = 0.2941585663520804
$ws.Range("T9").Value = 0.3239463845959386

$ws.Range("I10").Value = 0.822180234441485
$ws.Range("J10").Value = 0.8526739017519405
$ws.Range("M10").Value = 82.43477250000001
$ws.Range("N10").Value = 164.869545
$ws.Range("O10").Value = 0.1748236883957081
$ws.Range("P10").Value = 0.1237612588479007
$ws.Range("Q10").Value = 978.1391357063002
$ws.Range("R10").Value = 5868.834814237801
$ws.Range("S10").Value = 0.1437365811111084
$ws.Range("T10").Value = 0.1055279954675713

$ws.Range("I11").Value = 0.822180234441485
$ws.Range("J11").Value = 0.8526739017519405
$ws.Range("M11").Value = 123.1883796666667
$ws.Range("N11").Value = 369.565139
$ws.Range("O11").Value = 0.2612519722889696
$ws.Range("P11").Value = 0.2774184087724594
$ws.Range("Q11").Value = 1461.705680284529
$ws.Range("R11").Value = 13155.35112256076
$ws.Range("S11").Value = 0.2147962078248453
$ws.Range("T11").Value = 0.2365474370258277

$ws.Range("G12").Value = 0.37892
$ws.Range("H12").Value = 1.13676
$ws.Range("I12").Value = 0.02625574638939025
$ws.Range("J12").Value = 0.02722954016579943
$ws.Range("M12").Value = 16.14072933333334
$ws.Range("N12").Value = 48.42218800000001
$ws.Range("O12").Value = 0.03423048004954622
$ws.Range("P12").Value = 0.03634868370049611
$ws.Range("Q12").Value = 6.116045158986667
$ws.Range("R12").Value = 55.04440643088001
$ws.Range("S12").Value = 0.0008987468029679681
$ws.Range("T12").Value = 0.0009897579427965979

$ws.Range("G13").Value = 0.37892
$ws.Range("H13").Value = 1.13676
$ws.Range("I13").Value = 0.02625574638939025
$ws.Range("J13").Value = 0.02722954016579943
$ws.Range("O13").Value = 0.1719151703242873
$ws.Range("P13").Value = 0.1825533892714798
$ws.Range("Q13").Value = 30.71651182502666
$ws.Range("R13").Value = 276.44860642524
$ws.Range("S13").Value = 0.004513761112523317
$ws.Range("T13").Value = 0.004970844845570579

$ws.Range("G14").Value = 0.37892
$ws.Range("H14").Value = 1.13676
$ws.Range("I14").Value = 0.02625574638939025
$ws.Range("J14").Value = 0.02722954016579943
$ws.Range("M14").Value = 168.70371
$ws.Range("N14").Value = 506.11113
$ws.Range("O14").Value = 0.3577786889414888
$ws.Range("P14").Value = 0.3799182594076638
$ws.Range("Q14").Value = 63.9252097932
$ws.Range("R14").Value = 575.3268881388
$ws.Range("S14").Value = 0.009393746520376271
$ws.Range("T14").Value = 0.01034499950426159

$ws.Range("G15").Value = 0.37892
$ws.Range("H15").Value = 1.13676
$ws.Range("I15").Value = 0.02625574638939025
$ws.Range("J15").Value = 0.02722954016579943
$ws.Range("M15").Value = 82.43477250000001
$ws.Range("N15").Value = 164.869545
$ws.Range("O15").Value = 0.1748236883957081
$ws.Range("P15").Value = 0.1237612588479007
$ws.Range("Q15").Value = 31.2361839957
$ws.Range("R15").Value = 187.4171039742
$ws.Range("S15").Value = 0.004590126425375498
$ws.Range("T15").Value = 0.003369962168768812

$ws.Range("G16").Value = 0.37892
$ws.Range("H16").Value = 1.13676
$ws.Range("I16").Value = 0.02625574638939025
$ws.Range("J16").Value = 0.02722954016579943
$ws.Range("M16").Value = 123.1883796666667
$ws.Range("N16").Value = 369.565139
$ws.Range("O16").Value = 0.2612519722889696
$ws.Range("P16").Value = 0.2774184087724594
$ws.Range("Q16").Value = 46.67854082329333
$ws.Range("R16").Value = 420.10686740964
$ws.Range("S16").Value = 0.006859365528147194
$ws.Range("T16").Value = 0.007553975704401849

$ws.Range("G17").Value = 1.548357
$ws.Range("H17").Value = 3.096714
$ws.Range("I17").Value = 0.1072872076222874
$ws.Range("J17").Value = 0.0741775733180209
$ws.Range("M17").Value = 16.14072933333334
$ws.Range("N17").Value = 48.42218800000001
$ws.Range("O17").Value = 0.03423048004954622
$ws.Range("P17").Value = 0.03634868370049611
$ws.Range("Q17").Value = 24.99161124837201
$ws.Range("R17").Value = 149.949667490232
$ws.Range("S17").Value = 0.00367249262008623
$ws.Range("T17").Value = 0.002696257150207102

$ws.Range("G18").Value = 1.548357
$ws.Range("H18").Value = 3.096714
$ws.Range("I18").Value = 0.1072872076222874
$ws.Range("J18").Value = 0.0741775733180209
$ws.Range("O18").Value = 0.1719151703242873
$ws.Range("P18").Value = 0.1825533892714798
$ws.Range("Q18").Value = 125.514953287931
$ws.Range("R18").Value = 753.0897197275862
$ws.Range("S18").Value = 0.01844429857200271
$ws.Range("T18").Value = 0.01354136741713841

$ws.Range("G19").Value = 1.548357
$ws.Range("H19").Value = 3.096714
$ws.Range("I19").Value = 0.1072872076222874
$ws.Range("J19").Value = 0.0741775733180209
$ws.Range("M19").Value = 168.70371
$ws.Range("N19").Value = 506.11113
$ws.Range("O19").Value = 0.3577786889414888
$ws.Range("P19").Value = 0.3799182594076638
$ws.Range("Q19").Value = 261.21357030447
$ws.Range("R19").Value = 1567.28142182682
$ws.Range("S19").Value = 0.03838507648329528
$ws.Range("T19").Value = 0.02818141454206687

$ws.Range("G20").Value = 1.548357
$ws.Range("H20").Value = 3.096714
$ws.Range("I20").Value = 0.1072872076222874
$ws.Range("J20").Value = 0.0741775733180209
$ws.Range("M20").Value = 82.43477250000001
$ws.Range("N20").Value = 164.869545
$ws.Range("O20").Value = 0.1748236883957081
$ws.Range("P20").Value = 0.1237612588479007
$ws.Range("Q20").Value = 127.6384570437825
$ws.Range("R20").Value = 510.5538281751301
$ws.Range("S20").Value = 0.0187563453542044
$ws.Range("T20").Value = 0.009180309852120717

$ws.Range("G21").Value = 1.548357
$ws.Range("H21").Value = 3.096714
$ws.Range("I21").Value = 0.1072872076222874
$ws.Range("J21").Value = 0.0741775733180209
$ws.Range("M21").Value = 123.1883796666667
$ws.Range("N21").Value = 369.565139
$ws.Range("O21").Value = 0.2612519722889696
$ws.Range("P21").Value = 0.2774184087724594
$ws.Range("Q21").Value = 190.739589975541
$ws.Range("R21").Value = 1144.437539853246
$ws.Range("S21").Value = 0.02802899459269874
$ws.Range("T21").Value = 0.0205782243564878

$ws.Range("G22").Value = 0.16626
$ws.Range("H22").Value = 0.49878
$ws.Range("I22").Value = 0.01152032195371061
$ws.Range("J22").Value = 0.01194759671689489
$ws.Range("M22").Value = 16.14072933333334
$ws.Range("N22").Value = 48.42218800000001
$ws.Range("O22").Value = 0.03423048004954622
$ws.Range("P22").Value = 0.03634868370049611
$ws.Range("Q22").Value = 2.68355765896
$ws.Range("R22").Value = 24.15201893064
$ws.Range("S22").Value = 0.0003943461508008403
$ws.Range("T22").Value = 0.0004342794140434983

$ws.Range("G23").Value = 0.16626
$ws.Range("H23").Value = 0.49878
$ws.Range("I23").Value = 0.01152032195371061
$ws.Range("J23").Value = 0.01194759671689489
$ws.Range("O23").Value = 0.1719151703242873
$ws.Range("P23").Value = 0.1825533892714798
$ws.Range("Q23").Value = 13.47758697358
$ws.Range("R23").Value = 121.29828276222
$ws.Range("S23").Value = 0.001980518110862786
$ws.Range("T23").Value = 0.002181074274317968

$ws.Range("G24").Value = 0.16626
$ws.Range("H24").Value = 0.49878
$ws.Range("I24").Value = 0.01152032195371061
$ws.Range("J24").Value = 0.01194759671689489
$ws.Range("M24").Value = 168.70371
$ws.Range("N24").Value = 506.11113
$ws.Range("O24").Value = 0.3577786889414888
$ws.Range("P24").Value = 0.3799182594076638
$ws.Range("Q24").Value = 28.0486788246
$ws.Range("R24").Value = 252.4381094214
$ws.Range("S24").Value = 0.004121725684782432
$ws.Range("T24").Value = 0.004539110148787427

$ws.Range("G25").Value = 0.16626
$ws.Range("H25").Value = 0.49878
$ws.Range("I25").Value = 0.01152032195371061
$ws.Range("J25").Value = 0.01194759671689489
$ws.Range("M25").Value = 82.43477250000001
$ws.Range("N25").Value = 164.869545
$ws.Range("O25").Value = 0.1748236883957081
$ws.Range("P25").Value = 0.1237612588479007
$ws.Range("Q25").Value = 13.70560527585
$ws.Range("R25").Value = 82.23363165510001
$ws.Range("S25").Value = 0.002014025175453738
$ws.Range("T25").Value = 0.001478649609889957

$ws.Range("G26").Value = 0.16626
$ws.Range("H26").Value = 0.49878
$ws.Range("I26").Value = 0.01152032195371061
$ws.Range("J26").Value = 0.01194759671689489
$ws.Range("M26").Value = 123.1883796666667
$ws.Range("N26").Value = 369.565139
$ws.Range("O26").Value = 0.2612519722889696
$ws.Range("P26").Value = 0.2774184087724594
$ws.Range("Q26").Value = 20.48130000338
$ws.Range("R26").Value = 184.33170003042
$ws.Range("S26").Value = 0.003009706831810811
$ws.Range("T26").Value = 0.003314483269856042

